$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on price cells whose new values would otherwise
# be auto-parsed as numbers by Excel, losing the original text formatting
# (trailing zeros, punctuation) used by this coin-price list.
$textCells = @("D5", "D6", "D7", "D8", "D10", "D14", "D18", "D19", "D20", "D29", "D32", "D34", "D35", "D36", "D38", "D39", "D40", "D42", "D44", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "58.665.65"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").Value = "3.165.23"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "529.20"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "139.90"
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.536"
$ws.Range("E8").Value = "  +13.79%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "0.438"
$ws.Range("E10").Value = "  +6.70%  "
$ws.Range("E11").Value = "  +4.06%  "
$ws.Range("E12").Value = "  +2.54%  "
$ws.Range("D13").Value = "3.710.37"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").Value = "25.73"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("E15").Value = "  +3.51%  "
$ws.Range("D16").Value = "58.711.45"
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("D17").Value = "3.174.18"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "6.24"
$ws.Range("E18").Value = "  +3.41%  "
$ws.Range("D19").Value = "12.97"
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("D20").Value = "375.89"
$ws.Range("E20").Value = "  +4.31%  "
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  +4.86%  "
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("E27").Value = "  +14.13%  "
$ws.Range("D28").Value = "0.0₃0862"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").Value = "22.45"
$ws.Range("E29").Value = "  +5.33%  "
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("D32").Value = "5.14"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").Value = "6.31"
$ws.Range("E34").Value = "  +4.14%  "
$ws.Range("D35").Value = "156.78"
$ws.Range("E35").Value = "  -1.51%  "
$ws.Range("D36").Value = "1.33"
$ws.Range("E36").Value = "  +3.47%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "2.693.88"
$ws.Range("E37").Value = "  +8.38%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").Value = "25.02"
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "1.69"
$ws.Range("E39").Value = "  +1.52%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.0693"
$ws.Range("E40").Value = "  +3.90%  "
$ws.Range("E41").Value = "  +6.67%  "
$ws.Range("D42").Value = "0.723"
$ws.Range("E42").Value = "  +3.85%  "
$ws.Range("E43").Value = "  +7.55%  "
$ws.Range("D44").Value = "39.12"
$ws.Range("E44").Value = "  +3.70%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("E47").Value = "  +13.43%  "
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").Value = "20.04"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("D51").Value = "0.749"
$ws.Range("E51").Value = "  +1.29%  "
